# Updated symbol list on Fri Dec 16 07:56:36 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The "Price" column (D) stores numeric-looking values as literal TEXT
# (trailing zeros / exact digit counts matter, e.g. "0.1590", "1.350").
# Plain `.Value = "263.74"` assignment lets COM auto-convert the string
# to a real number (losing the formatting / exact text), so for every
# such cell we: switch the cell to text format, assign the literal
# string, then restore the "Normal" style so no stray formatting is
# left behind (matches the source workbook's unstyled data cells).
# ---------------------------------------------------------------------
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Simple price/volume refreshes (single-cell updates)
Set-TextValue "D2"  "263.74"
Set-TextValue "D3"  "23.19"
Set-TextValue "D4"  "6.189"
Set-TextValue "D5"  "0.06248"
Set-TextValue "D6"  "6.738"
Set-TextValue "D8"  "1.350"
Set-TextValue "D9"  "0.7966"

# Rows 10-18: coin list shifted up by one rank, new values/prices
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1590"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.08184"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03414"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03085"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09340"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "3.711"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001689"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04795"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0006132"
$ws.Range("E18").Value = "17OneONEWorstin24h"

# More single-cell price refreshes
Set-TextValue "D20" "0.006185"
Set-TextValue "D21" "0.001098"
Set-TextValue "D22" "0.0001497"
Set-TextValue "D23" "3.702"
Set-TextValue "D24" "2.198"
Set-TextValue "D25" "0.3340"
Set-TextValue "D26" "0.1273"
Set-TextValue "D27" "0.0003199"
Set-TextValue "D40" "0.04633"

# Rows 41-43: coin list shifted up by one rank
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1124"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003144"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003358"
$ws.Range("E43").Value = "42KickTokenKICK"

# Remaining single-cell refreshes
Set-TextValue "D46" "0.00005889"
Set-TextValue "D48" "0.6993"

Set-TextValue "D49" "0.1339"
$ws.Range("E49").Value = "48BOLOBOLO"

Set-TextValue "D50" "0.00002098"
